# Move the "_GoBack" bookmark from the empty paragraph that follows the
# "Patient Gender / Patient Dosing Weight" table up into the weight-
# placeholder cell, where it replaces the trailing run of four spaces.

$d = $word.ActiveDocument

# Locate the table holding "Patient Dosing Weight (kg):" -> the last cell
# in its single row contains "{weight_placeholder}    " (with 4 trailing
# spaces in their own run).
$targetTable = $null
foreach ($tbl in $d.Tables) {
    if ($tbl.Cell(1, 1).Range.Text -like "*Patient Gender*") {
        $targetTable = $tbl
        break
    }
}

$cell = $targetTable.Cell(1, 4)
$cellRange = $cell.Range

# $cellRange.End points one past the cell-mark; back up 1 to land right
# after the real paragraph content (before the paragraph mark), then 4
# more to land right before the run of 4 trailing spaces.
$contentEnd = $cellRange.End - 1
$spacesStart = $contentEnd - 4

# Create (or relocate, since the name already exists elsewhere in the
# document) the _GoBack bookmark as a zero-length range right after "}"
# and before the spaces, then delete the now-orphaned space run. Adding
# the bookmark before deleting keeps its anchor position stable.
$bmRange = $d.Range($spacesStart, $spacesStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

$spacesRange = $d.Range($spacesStart, $contentEnd)
$spacesRange.Delete()
